$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 819 (A819=818)
$ws.Cells.Item(819, 1).Value = 818
$ws.Cells.Item(819, 2).Value = 100
$ws.Cells.Item(819, 3).Value = 43220
$ws.Cells.Item(819, 4).Value = "April"
$ws.Cells.Item(819, 5).Value = 2018
$ws.Cells.Item(819, 6).Value = "Monday"
$ws.Cells.Item(819, 7).Value = "Bicep Curl"
$ws.Cells.Item(819, 8).Value = 40
$ws.Cells.Item(819, 9).Value = 4
$ws.Cells.Item(819, 10).Value = 8
$ws.Cells.Item(819, 11).Value = "Arms"

# Row 820 (A820=819)
$ws.Cells.Item(820, 1).Value = 819
$ws.Cells.Item(820, 2).Value = 100
$ws.Cells.Item(820, 3).Value = 43220
$ws.Cells.Item(820, 4).Value = "April"
$ws.Cells.Item(820, 5).Value = 2018
$ws.Cells.Item(820, 6).Value = "Monday"
$ws.Cells.Item(820, 7).Value = "Laying down tricep curl"
$ws.Cells.Item(820, 8).Value = 30
$ws.Cells.Item(820, 9).Value = 4
$ws.Cells.Item(820, 10).Value = 8
$ws.Cells.Item(820, 11).Value = "Arms"

# Row 821 (A821=820)
$ws.Cells.Item(821, 1).Value = 820
$ws.Cells.Item(821, 2).Value = 100
$ws.Cells.Item(821, 3).Value = 43220
$ws.Cells.Item(821, 4).Value = "April"
$ws.Cells.Item(821, 5).Value = 2018
$ws.Cells.Item(821, 6).Value = "Monday"
$ws.Cells.Item(821, 7).Value = "Hammer Curl"
$ws.Cells.Item(821, 8).Value = 20
$ws.Cells.Item(821, 9).Value = 4
$ws.Cells.Item(821, 10).Value = 8
$ws.Cells.Item(821, 11).Value = "Arms"

# Row 822 (A822=821)
$ws.Cells.Item(822, 1).Value = 821
$ws.Cells.Item(822, 2).Value = 100
$ws.Cells.Item(822, 3).Value = 43220
$ws.Cells.Item(822, 4).Value = "April"
$ws.Cells.Item(822, 5).Value = 2018
$ws.Cells.Item(822, 6).Value = "Monday"
$ws.Cells.Item(822, 7).Value = "Tricep Pull down"
$ws.Cells.Item(822, 8).Value = 45
$ws.Cells.Item(822, 9).Value = 4
$ws.Cells.Item(822, 10).Value = 8
$ws.Cells.Item(822, 11).Value = "Arms"

# Row 823 (A823=822)
$ws.Cells.Item(823, 1).Value = 822
$ws.Cells.Item(823, 2).Value = 100
$ws.Cells.Item(823, 3).Value = 43220
$ws.Cells.Item(823, 4).Value = "April"
$ws.Cells.Item(823, 5).Value = 2018
$ws.Cells.Item(823, 6).Value = "Monday"
$ws.Cells.Item(823, 7).Value = "Plank"
$ws.Cells.Item(823, 8).Value = 0
$ws.Cells.Item(823, 9).Value = 3
$ws.Cells.Item(823, 10).Value = 30
$ws.Cells.Item(823, 11).Value = "Core"

# Row 824 (A824=823)
$ws.Cells.Item(824, 1).Value = 823
$ws.Cells.Item(824, 2).Value = 100
$ws.Cells.Item(824, 3).Value = 43220
$ws.Cells.Item(824, 4).Value = "April"
$ws.Cells.Item(824, 5).Value = 2018
$ws.Cells.Item(824, 6).Value = "Monday"
$ws.Cells.Item(824, 7).Value = "Left Plank"
$ws.Cells.Item(824, 8).Value = 0
$ws.Cells.Item(824, 9).Value = 3
$ws.Cells.Item(824, 10).Value = 30
$ws.Cells.Item(824, 11).Value = "Core"

# Row 825 (A825=824)
$ws.Cells.Item(825, 1).Value = 824
$ws.Cells.Item(825, 2).Value = 100
$ws.Cells.Item(825, 3).Value = 43220
$ws.Cells.Item(825, 4).Value = "April"
$ws.Cells.Item(825, 5).Value = 2018
$ws.Cells.Item(825, 6).Value = "Monday"
$ws.Cells.Item(825, 7).Value = "Right Plank"
$ws.Cells.Item(825, 8).Value = 0
$ws.Cells.Item(825, 9).Value = 3
$ws.Cells.Item(825, 10).Value = 30
$ws.Cells.Item(825, 11).Value = "Core"

# Row 826 (A826=825)
$ws.Cells.Item(826, 1).Value = 825
$ws.Cells.Item(826, 2).Value = 101
$ws.Cells.Item(826, 3).Value = 43222
$ws.Cells.Item(826, 4).Value = "May"
$ws.Cells.Item(826, 5).Value = 2018
$ws.Cells.Item(826, 6).Value = "Wednesday"
$ws.Cells.Item(826, 7).Value = "Bench Press"
$ws.Cells.Item(826, 8).Value = 85
$ws.Cells.Item(826, 9).Value = 5
$ws.Cells.Item(826, 10).Value = 5
$ws.Cells.Item(826, 11).Value = "Chest"

# Row 827 (A827=826)
$ws.Cells.Item(827, 1).Value = 826
$ws.Cells.Item(827, 2).Value = 101
$ws.Cells.Item(827, 3).Value = 43222
$ws.Cells.Item(827, 4).Value = "May"
$ws.Cells.Item(827, 5).Value = 2018
$ws.Cells.Item(827, 6).Value = "Wednesday"
$ws.Cells.Item(827, 7).Value = "Overhead Press"
$ws.Cells.Item(827, 8).Value = 52.5
$ws.Cells.Item(827, 9).Value = 5
$ws.Cells.Item(827, 10).Value = 5
$ws.Cells.Item(827, 11).Value = "Shoulders"

# Row 828 (A828=827)
$ws.Cells.Item(828, 1).Value = 827
$ws.Cells.Item(828, 2).Value = 101
$ws.Cells.Item(828, 3).Value = 43222
$ws.Cells.Item(828, 4).Value = "May"
$ws.Cells.Item(828, 5).Value = 2018
$ws.Cells.Item(828, 6).Value = "Wednesday"
$ws.Cells.Item(828, 7).Value = "Laterial Raises"
$ws.Cells.Item(828, 8).Value = 10
$ws.Cells.Item(828, 9).Value = 4
$ws.Cells.Item(828, 10).Value = 8
$ws.Cells.Item(828, 11).Value = "Shoulders"

# Row 829 (A829=828)
$ws.Cells.Item(829, 1).Value = 828
$ws.Cells.Item(829, 2).Value = 101
$ws.Cells.Item(829, 3).Value = 43222
$ws.Cells.Item(829, 4).Value = "May"
$ws.Cells.Item(829, 5).Value = 2018
$ws.Cells.Item(829, 6).Value = "Wednesday"
$ws.Cells.Item(829, 7).Value = "Front raises"
$ws.Cells.Item(829, 8).Value = 10
$ws.Cells.Item(829, 9).Value = 4
$ws.Cells.Item(829, 10).Value = 8
$ws.Cells.Item(829, 11).Value = "Shoulders"

# Row 830 (A830=829)
$ws.Cells.Item(830, 1).Value = 829
$ws.Cells.Item(830, 2).Value = 101
$ws.Cells.Item(830, 3).Value = 43222
$ws.Cells.Item(830, 4).Value = "May"
$ws.Cells.Item(830, 5).Value = 2018
$ws.Cells.Item(830, 6).Value = "Wednesday"
$ws.Cells.Item(830, 7).Value = "Upright Rows"
$ws.Cells.Item(830, 8).Value = 30
$ws.Cells.Item(830, 9).Value = 4
$ws.Cells.Item(830, 10).Value = 12
$ws.Cells.Item(830, 11).Value = "Shoulders"

# Row 831 (A831=830)
$ws.Cells.Item(831, 1).Value = 830
$ws.Cells.Item(831, 2).Value = 102
$ws.Cells.Item(831, 3).Value = 43224
$ws.Cells.Item(831, 4).Value = "May"
$ws.Cells.Item(831, 5).Value = 2018
$ws.Cells.Item(831, 6).Value = "Friday"
$ws.Cells.Item(831, 7).Value = "Pec Fly"
$ws.Cells.Item(831, 8).Value = 105
$ws.Cells.Item(831, 9).Value = 4
$ws.Cells.Item(831, 10).Value = 8
$ws.Cells.Item(831, 11).Value = "Chest"

# Row 832 (A832=831)
$ws.Cells.Item(832, 1).Value = 831
$ws.Cells.Item(832, 2).Value = 102
$ws.Cells.Item(832, 3).Value = 43224
$ws.Cells.Item(832, 4).Value = "May"
$ws.Cells.Item(832, 5).Value = 2018
$ws.Cells.Item(832, 6).Value = "Friday"
$ws.Cells.Item(832, 7).Value = "Incline Bench"
$ws.Cells.Item(832, 8).Value = 70
$ws.Cells.Item(832, 9).Value = 3
$ws.Cells.Item(832, 10).Value = 8
$ws.Cells.Item(832, 11).Value = "Chest"

# Row 833 (A833=832)
$ws.Cells.Item(833, 1).Value = 832
$ws.Cells.Item(833, 2).Value = 102
$ws.Cells.Item(833, 3).Value = 43224
$ws.Cells.Item(833, 4).Value = "May"
$ws.Cells.Item(833, 5).Value = 2018
$ws.Cells.Item(833, 6).Value = "Friday"
$ws.Cells.Item(833, 7).Value = "Seated Row"
$ws.Cells.Item(833, 8).Value = 70
$ws.Cells.Item(833, 9).Value = 4
$ws.Cells.Item(833, 10).Value = 8
$ws.Cells.Item(833, 11).Value = "Back"

# Row 834 (A834=833)
$ws.Cells.Item(834, 1).Value = 833
$ws.Cells.Item(834, 2).Value = 102
$ws.Cells.Item(834, 3).Value = 43224
$ws.Cells.Item(834, 4).Value = "May"
$ws.Cells.Item(834, 5).Value = 2018
$ws.Cells.Item(834, 6).Value = "Friday"
$ws.Cells.Item(834, 7).Value = "Dumbell Rows"
$ws.Cells.Item(834, 8).Value = 30
$ws.Cells.Item(834, 9).Value = 4
$ws.Cells.Item(834, 10).Value = 8
$ws.Cells.Item(834, 11).Value = "Back"

# Row 835 (A835=834)
$ws.Cells.Item(835, 1).Value = 834
$ws.Cells.Item(835, 2).Value = 102
$ws.Cells.Item(835, 3).Value = 43224
$ws.Cells.Item(835, 4).Value = "May"
$ws.Cells.Item(835, 5).Value = 2018
$ws.Cells.Item(835, 6).Value = "Friday"
$ws.Cells.Item(835, 7).Value = "Heel-taps"
$ws.Cells.Item(835, 8).Value = 0
$ws.Cells.Item(835, 9).Value = 2
$ws.Cells.Item(835, 10).Value = 10
$ws.Cells.Item(835, 11).Value = "Core"

# Row 836 (A836=835)
$ws.Cells.Item(836, 1).Value = 835
$ws.Cells.Item(836, 2).Value = 102
$ws.Cells.Item(836, 3).Value = 43224
$ws.Cells.Item(836, 4).Value = "May"
$ws.Cells.Item(836, 5).Value = 2018
$ws.Cells.Item(836, 6).Value = "Friday"
$ws.Cells.Item(836, 7).Value = "Leg Raises"
$ws.Cells.Item(836, 8).Value = 0
$ws.Cells.Item(836, 9).Value = 2
$ws.Cells.Item(836, 10).Value = 10
$ws.Cells.Item(836, 11).Value = "Core"

# Row 837 (A837=836)
$ws.Cells.Item(837, 1).Value = 836
$ws.Cells.Item(837, 2).Value = 102
$ws.Cells.Item(837, 3).Value = 43224
$ws.Cells.Item(837, 4).Value = "May"
$ws.Cells.Item(837, 5).Value = 2018
$ws.Cells.Item(837, 6).Value = "Friday"
$ws.Cells.Item(837, 7).Value = "Scissors"
$ws.Cells.Item(837, 8).Value = 0
$ws.Cells.Item(837, 9).Value = 2
$ws.Cells.Item(837, 10).Value = 12
$ws.Cells.Item(837, 11).Value = "Core"

# Row 838 (A838=837)
$ws.Cells.Item(838, 1).Value = 837
$ws.Cells.Item(838, 2).Value = 102
$ws.Cells.Item(838, 3).Value = 43224
$ws.Cells.Item(838, 4).Value = "May"
$ws.Cells.Item(838, 5).Value = 2018
$ws.Cells.Item(838, 6).Value = "Friday"
$ws.Cells.Item(838, 7).Value = "Knee-Pull ins"
$ws.Cells.Item(838, 8).Value = 0
$ws.Cells.Item(838, 9).Value = 2
$ws.Cells.Item(838, 10).Value = 10
$ws.Cells.Item(838, 11).Value = "Core"

# Row 839 (A839=838)
$ws.Cells.Item(839, 1).Value = 838
$ws.Cells.Item(839, 2).Value = 102
$ws.Cells.Item(839, 3).Value = 43224
$ws.Cells.Item(839, 4).Value = "May"
$ws.Cells.Item(839, 5).Value = 2018
$ws.Cells.Item(839, 6).Value = "Friday"
$ws.Cells.Item(839, 7).Value = "Flitter Kicks"
$ws.Cells.Item(839, 8).Value = 0
$ws.Cells.Item(839, 9).Value = 2
$ws.Cells.Item(839, 10).Value = 10
$ws.Cells.Item(839, 11).Value = "Core"

# Row 840 (A840=839)
$ws.Cells.Item(840, 1).Value = 839
$ws.Cells.Item(840, 2).Value = 103
$ws.Cells.Item(840, 3).Value = 43225
$ws.Cells.Item(840, 4).Value = "May"
$ws.Cells.Item(840, 5).Value = 2018
$ws.Cells.Item(840, 6).Value = "Saturday"
$ws.Cells.Item(840, 7).Value = "Shoulder Press"
$ws.Cells.Item(840, 8).Value = 25
$ws.Cells.Item(840, 9).Value = 4
$ws.Cells.Item(840, 10).Value = 8
$ws.Cells.Item(840, 11).Value = "Shoulders"

# Row 841 (A841=840)
$ws.Cells.Item(841, 1).Value = 840
$ws.Cells.Item(841, 2).Value = 103
$ws.Cells.Item(841, 3).Value = 43225
$ws.Cells.Item(841, 4).Value = "May"
$ws.Cells.Item(841, 5).Value = 2018
$ws.Cells.Item(841, 6).Value = "Saturday"
$ws.Cells.Item(841, 7).Value = "Shoulder Shrug"
$ws.Cells.Item(841, 8).Value = 25
$ws.Cells.Item(841, 9).Value = 4
$ws.Cells.Item(841, 10).Value = 8
$ws.Cells.Item(841, 11).Value = "Shoulders"

# Row 842 (A842=841)
$ws.Cells.Item(842, 1).Value = 841
$ws.Cells.Item(842, 2).Value = 103
$ws.Cells.Item(842, 3).Value = 43225
$ws.Cells.Item(842, 4).Value = "May"
$ws.Cells.Item(842, 5).Value = 2018
$ws.Cells.Item(842, 6).Value = "Saturday"
$ws.Cells.Item(842, 7).Value = "Tricep Pull down"
$ws.Cells.Item(842, 8).Value = 45
$ws.Cells.Item(842, 9).Value = 4
$ws.Cells.Item(842, 10).Value = 8
$ws.Cells.Item(842, 11).Value = "Arms"

# Row 843 (A843=842)
$ws.Cells.Item(843, 1).Value = 842
$ws.Cells.Item(843, 2).Value = 103
$ws.Cells.Item(843, 3).Value = 43225
$ws.Cells.Item(843, 4).Value = "May"
$ws.Cells.Item(843, 5).Value = 2018
$ws.Cells.Item(843, 6).Value = "Saturday"
$ws.Cells.Item(843, 7).Value = "Hammer Curl"
$ws.Cells.Item(843, 8).Value = 20
$ws.Cells.Item(843, 9).Value = 4
$ws.Cells.Item(843, 10).Value = 8
$ws.Cells.Item(843, 11).Value = "Arms"

# Row 844 (A844=843)
$ws.Cells.Item(844, 1).Value = 843
$ws.Cells.Item(844, 2).Value = 103
$ws.Cells.Item(844, 3).Value = 43225
$ws.Cells.Item(844, 4).Value = "May"
$ws.Cells.Item(844, 5).Value = 2018
$ws.Cells.Item(844, 6).Value = "Saturday"
$ws.Cells.Item(844, 7).Value = "Russian Twists"
$ws.Cells.Item(844, 8).Value = 10
$ws.Cells.Item(844, 9).Value = 4
$ws.Cells.Item(844, 10).Value = 12
$ws.Cells.Item(844, 11).Value = "Core"

# Row 845 (A845=844)
$ws.Cells.Item(845, 1).Value = 844
$ws.Cells.Item(845, 2).Value = 103
$ws.Cells.Item(845, 3).Value = 43225
$ws.Cells.Item(845, 4).Value = "May"
$ws.Cells.Item(845, 5).Value = 2018
$ws.Cells.Item(845, 6).Value = "Saturday"
$ws.Cells.Item(845, 7).Value = "Left Situp"
$ws.Cells.Item(845, 8).Value = 0
$ws.Cells.Item(845, 9).Value = 4
$ws.Cells.Item(845, 10).Value = 12
$ws.Cells.Item(845, 11).Value = "Core"

# Row 846 (A846=845)
$ws.Cells.Item(846, 1).Value = 845
$ws.Cells.Item(846, 2).Value = 103
$ws.Cells.Item(846, 3).Value = 43225
$ws.Cells.Item(846, 4).Value = "May"
$ws.Cells.Item(846, 5).Value = 2018
$ws.Cells.Item(846, 6).Value = "Saturday"
$ws.Cells.Item(846, 7).Value = "Right Situp"
$ws.Cells.Item(846, 8).Value = 0
$ws.Cells.Item(846, 9).Value = 4
$ws.Cells.Item(846, 10).Value = 12
$ws.Cells.Item(846, 11).Value = "Core"

# Row 847 (A847=846)
$ws.Cells.Item(847, 1).Value = 846
$ws.Cells.Item(847, 2).Value = 104
$ws.Cells.Item(847, 3).Value = 43226
$ws.Cells.Item(847, 4).Value = "May"
$ws.Cells.Item(847, 5).Value = 2018
$ws.Cells.Item(847, 6).Value = "Sunday"
$ws.Cells.Item(847, 7).Value = "Barbell Squat"
$ws.Cells.Item(847, 8).Value = 75
$ws.Cells.Item(847, 9).Value = 3
$ws.Cells.Item(847, 10).Value = 8
$ws.Cells.Item(847, 11).Value = "Legs"

# Row 848 (A848=847)
$ws.Cells.Item(848, 1).Value = 847
$ws.Cells.Item(848, 2).Value = 104
$ws.Cells.Item(848, 3).Value = 43226
$ws.Cells.Item(848, 4).Value = "May"
$ws.Cells.Item(848, 5).Value = 2018
$ws.Cells.Item(848, 6).Value = "Sunday"
$ws.Cells.Item(848, 7).Value = "Barbell Lunge"
$ws.Cells.Item(848, 8).Value = 45
$ws.Cells.Item(848, 9).Value = 3
$ws.Cells.Item(848, 10).Value = 8
$ws.Cells.Item(848, 11).Value = "Legs"

# Row 849 (A849=848)
$ws.Cells.Item(849, 1).Value = 848
$ws.Cells.Item(849, 2).Value = 104
$ws.Cells.Item(849, 3).Value = 43226
$ws.Cells.Item(849, 4).Value = "May"
$ws.Cells.Item(849, 5).Value = 2018
$ws.Cells.Item(849, 6).Value = "Sunday"
$ws.Cells.Item(849, 7).Value = "Hip adduction"
$ws.Cells.Item(849, 8).Value = 65
$ws.Cells.Item(849, 9).Value = 3
$ws.Cells.Item(849, 10).Value = 12
$ws.Cells.Item(849, 11).Value = "Legs"

# Row 850 (A850=849)
$ws.Cells.Item(850, 1).Value = 849
$ws.Cells.Item(850, 2).Value = 104
$ws.Cells.Item(850, 3).Value = 43226
$ws.Cells.Item(850, 4).Value = "May"
$ws.Cells.Item(850, 5).Value = 2018
$ws.Cells.Item(850, 6).Value = "Sunday"
$ws.Cells.Item(850, 7).Value = "Hip abduction"
$ws.Cells.Item(850, 8).Value = 65
$ws.Cells.Item(850, 9).Value = 3
$ws.Cells.Item(850, 10).Value = 12
$ws.Cells.Item(850, 11).Value = "Legs"

# Row 851 (A851=850)
$ws.Cells.Item(851, 1).Value = 850
$ws.Cells.Item(851, 2).Value = 104
$ws.Cells.Item(851, 3).Value = 43226
$ws.Cells.Item(851, 4).Value = "May"
$ws.Cells.Item(851, 5).Value = 2018
$ws.Cells.Item(851, 6).Value = "Sunday"
$ws.Cells.Item(851, 7).Value = "Leg Extension"
$ws.Cells.Item(851, 8).Value = 108
$ws.Cells.Item(851, 9).Value = 4
$ws.Cells.Item(851, 10).Value = 8
$ws.Cells.Item(851, 11).Value = "Legs"

$ws.Range("A852").Select()